# finnegan_kyle.xlsx save_data regen:
#   - "Strike#" column (G) renumbered/recalculated as "K" (strikeouts)
#   - row 40's IP (H) and I0 (I) also corrected as part of the same regen
#
# Apply the new literal values cell-by-cell on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "G2"  = 1
    "G3"  = 2
    "G4"  = 1
    "G5"  = 1
    "G6"  = 1
    "G7"  = 1
    "G8"  = 0
    "G9"  = 1
    "G10" = 1
    "G11" = 2
    "G12" = 2
    "G13" = 4
    "G14" = 0
    "G16" = 2
    "G17" = 1
    "G18" = 1
    "G19" = 0
    "G20" = 0
    "G22" = 0
    "G23" = 1
    "G24" = 0
    "G25" = 0
    "G26" = 2
    "G27" = 0
    "G28" = 2
    "G29" = 2
    "G30" = 0
    "G31" = 0
    "G32" = 1
    "G33" = 3
    "G34" = 0
    "G35" = 0
    "G36" = 1
    "G37" = 1
    "G38" = 1
    "G39" = 2
    "G40" = 1
    "H40" = 3
    "I40" = 6
    "G41" = 2
    "G42" = 1
    "G43" = 2
    "G44" = 2
    "G45" = 1
    "G46" = 0
    "G47" = 1
    "G48" = 1
    "G49" = 0
    "G50" = 1
    "G51" = 0
    "G52" = 0
    "G53" = 1
    "G54" = 1
    "G55" = 3
    "G56" = 1
    "G57" = 1
    "G58" = 1
    "G59" = 1
    "G60" = 1
    "G61" = 2
    "G62" = 1
    "G63" = 0
    "G64" = 2
    "G65" = 1
    "G66" = 3
    "G67" = 0
    "G68" = 2
    "G69" = 1
    "G70" = 2
    "G72" = 1
    "G73" = 1
    "G74" = 0
    "G75" = 1
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
